# Seminar 1 - remarks added after the lab on 18.10.2021
#
# 1) The "_GoBack" bookmark (left over from the author's last edit before
#    saving in real Word) moves from its old spot near the picture to the
#    paragraph that starts with "Минимум действий..." - it now wraps
#    everything up to (but not including) the final full stop of that
#    paragraph, leaving the final "." in its own run.
# 2) The following paragraph ("Таблица глобальных дескрипторов содержит
#    дескрипторы сегментов физической памяти.") is made bold, and a
#    second bookmark ("_Hlk85462764") is added around the word "содержит ".

$d = $word.ActiveDocument

# --- Step 1: drop the stale _GoBack bookmark wherever it currently sits ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: locate the "Минимум действий..." paragraph -------------------
$hit1 = $d.Content.Duplicate
$hit1.Find.Execute("Минимум действий", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
$p1 = $hit1.Paragraphs(1).Range

$p1Text = $p1.Text
$periodPos = $p1Text.LastIndexOf(".")

# range from the start of the paragraph up to (but excluding) the final "."
$goBackRange = $d.Range($p1.Start, $p1.Start + $periodPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- Step 3: the "Таблица глобальных дескрипторов..." paragraph ----------
$hit2 = $d.Content.Duplicate
$hit2.Find.Execute("Таблица глобальных дескрипторов", $true, $false, $false,
                    $false, $false, $true, 1, $false, "", 0)
$p2 = $hit2.Paragraphs(1).Range

# whole paragraph, including the paragraph mark, becomes bold
$p2.Bold = 1

$word1 = "Таблица глобальных дескрипторов "
$word2 = "содержит "
$hlkStart = $p2.Start + $word1.Length
$hlkEnd = $hlkStart + $word2.Length

$hlkRange = $d.Range($hlkStart, $hlkEnd)
$d.Bookmarks.Add("_Hlk85462764", $hlkRange)
